# Update Data Sources from LFX: re-point table style references from the
# old table style GUID to the new one on every slide that uses it.

$oldStyleId = "{74CF3CD3-5AD3-4896-9624-B6345D4F0AB6}"
$newStyleId = "{C320EC70-FFC7-4695-BD8A-B35883984C6E}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
